# Weekly fruit/vegetable price update: Hortaliza - Terminal Hortofruticola Agro Chillan - Alcachofa
# Reassigns the Fecha/Variedad/Calidad/Volumen/Precio/Origen values across existing data rows
# (rows 2-35) to reflect the corrected weekly price records.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44488
$ws.Range("J2").Value = 100

$ws.Range("D3").Value = 44498
$ws.Range("J3").Value = 60
$ws.Range("K3").Value = 10500
$ws.Range("L3").Value = 11000
$ws.Range("M3").Value = 10750
$ws.Range("P3").Value = 269

$ws.Range("D4").Value = 44491
$ws.Range("J4").Value = 100
$ws.Range("K4").Value = 11000
$ws.Range("L4").Value = 12000
$ws.Range("M4").Value = 11500
$ws.Range("P4").Value = 288

$ws.Range("D5").Value = 44399
$ws.Range("H5").Value = "Española"
$ws.Range("I5").Value = "Segunda"
$ws.Range("K5").Value = 15500
$ws.Range("L5").Value = 16000
$ws.Range("M5").Value = 15750
$ws.Range("O5").Value = "Provincia del Elquí"
$ws.Range("P5").Value = 394

$ws.Range("D6").Value = 44475

$ws.Range("D7").Value = 44489

$ws.Range("D8").Value = 44426
$ws.Range("J8").Value = 120
$ws.Range("K8").Value = 13000
$ws.Range("L8").Value = 14000
$ws.Range("M8").Value = 13500
$ws.Range("O8").Value = "Región del Maule"
$ws.Range("P8").Value = 338

$ws.Range("D9").Value = 44455
$ws.Range("J9").Value = 100
$ws.Range("K9").Value = 13000
$ws.Range("L9").Value = 14000
$ws.Range("M9").Value = 13500
$ws.Range("P9").Value = 338

$ws.Range("D10").Value = 44510
$ws.Range("K10").Value = 11000
$ws.Range("L10").Value = 12000
$ws.Range("M10").Value = 11500
$ws.Range("O10").Value = "Provincia del Elquí"
$ws.Range("P10").Value = 288

$ws.Range("D11").Value = 44484

$ws.Range("D12").Value = 44482
$ws.Range("J12").Value = 120
$ws.Range("K12").Value = 11000
$ws.Range("L12").Value = 12000
$ws.Range("M12").Value = 11500
$ws.Range("P12").Value = 288

$ws.Range("D13").Value = 44473
$ws.Range("J13").Value = 160

$ws.Range("D14").Value = 44515
$ws.Range("J14").Value = 120
$ws.Range("O14").Value = "Provincia del Elquí"

$ws.Range("D15").Value = 44516
$ws.Range("K15").Value = 11000
$ws.Range("L15").Value = 12000
$ws.Range("M15").Value = 11500
$ws.Range("P15").Value = 288

$ws.Range("D16").Value = 44446
$ws.Range("J16").Value = 160
$ws.Range("K16").Value = 12500
$ws.Range("L16").Value = 13000
$ws.Range("M16").Value = 12750
$ws.Range("P16").Value = 319

$ws.Range("D17").Value = 44503
$ws.Range("J17").Value = 160

$ws.Range("D18").Value = 44468
$ws.Range("J18").Value = 60
$ws.Range("K18").Value = 12000
$ws.Range("L18").Value = 13000
$ws.Range("M18").Value = 12500
$ws.Range("P18").Value = 312

$ws.Range("D19").Value = 44427
$ws.Range("J19").Value = 120
$ws.Range("K19").Value = 13000
$ws.Range("L19").Value = 14000
$ws.Range("M19").Value = 13500
$ws.Range("P19").Value = 338

$ws.Range("D20").Value = 44505
$ws.Range("J20").Value = 120
$ws.Range("K20").Value = 11000
$ws.Range("L20").Value = 12000
$ws.Range("M20").Value = 11500
$ws.Range("P20").Value = 288

$ws.Range("D21").Value = 44425
$ws.Range("J21").Value = 120
$ws.Range("K21").Value = 14000
$ws.Range("L21").Value = 15000
$ws.Range("M21").Value = 14500
$ws.Range("O21").Value = "Región del Maule"
$ws.Range("P21").Value = 362

$ws.Range("D22").Value = 44432
$ws.Range("J22").Value = 120
$ws.Range("K22").Value = 14000
$ws.Range("L22").Value = 15000
$ws.Range("M22").Value = 14500
$ws.Range("P22").Value = 362

$ws.Range("D23").Value = 44467
$ws.Range("J23").Value = 160
$ws.Range("O23").Value = "Provincia de Limarí"

$ws.Range("D24").Value = 44420
$ws.Range("H24").Value = "Madrigal"
$ws.Range("I24").Value = "Primera"
$ws.Range("K24").Value = 13000
$ws.Range("L24").Value = 14000
$ws.Range("M24").Value = 13500
$ws.Range("P24").Value = 338

$ws.Range("D25").Value = 44435
$ws.Range("K25").Value = 14000
$ws.Range("L25").Value = 15000
$ws.Range("M25").Value = 14500
$ws.Range("P25").Value = 362

$ws.Range("D26").Value = 44496

$ws.Range("D27").Value = 44508
$ws.Range("J27").Value = 160
$ws.Range("K27").Value = 11000
$ws.Range("L27").Value = 12000
$ws.Range("M27").Value = 11500
$ws.Range("P27").Value = 288

$ws.Range("D28").Value = 44512

$ws.Range("D29").Value = 44490

$ws.Range("D30").Value = 44417
$ws.Range("K30").Value = 15000
$ws.Range("L30").Value = 16000
$ws.Range("M30").Value = 15500
$ws.Range("P30").Value = 388

$ws.Range("D31").Value = 44495
$ws.Range("J31").Value = 120

$ws.Range("D32").Value = 44454

$ws.Range("D33").Value = 44494
$ws.Range("K33").Value = 11000
$ws.Range("L33").Value = 12000
$ws.Range("M33").Value = 11500
$ws.Range("P33").Value = 288

$ws.Range("D34").Value = 44453
$ws.Range("K34").Value = 12500
$ws.Range("L34").Value = 13000
$ws.Range("M34").Value = 12750
$ws.Range("P34").Value = 319

$ws.Range("D35").Value = 44487
$ws.Range("J35").Value = 100
